$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Onyeka Okongwu', 'PF,C', 'Atlanta Hawks'),
    @('Tari Eason', 'SF,PF', 'Houston Rockets'),
    @('Naz Reid', 'PF,C', 'Minnesota Timberwolves'),
    @('Julius Randle', 'PF,C', 'Minnesota Timberwolves'),
    @('Deandre Ayton', 'C', 'Portland Trail Blazers'),
    @('Collin Sexton', 'PG,SG', 'Utah Jazz'),
    @('Isaiah Hartenstein', 'C', 'Oklahoma City Thunder'),
    @('Damian Lillard', 'PG', 'Milwaukee Bucks'),
    @('Justin Edwards', 'SF', 'Philadelphia 76ers'),
    @('Derrick White', 'PG,SG', 'Boston Celtics'),
    @('Tyus Jones', 'PG', 'Phoenix Suns'),
    @('Coby White', 'PG,SG', 'Chicago Bulls'),
    @('Cade Cunningham', 'PG,SG', 'Detroit Pistons'),
    @('Devin Vassell', 'SG,SF', 'San Antonio Spurs'),
    @('Malik Monk', 'PG,SG,SF', 'Sacramento Kings'),
    @('Anthony Davis', 'PF,C', 'Los Angeles Lakers'),
    @('Cameron Johnson', 'SF,PF', 'Brooklyn Nets'),
    @('LaMelo Ball', 'PG,SG', 'Charlotte Hornets'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
